$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 199.09091
$ws.Range("I4").Value = 227.22223
$ws.Range("J4").Value = 72.5
$ws.Range("K4").Value = 227.22223
$ws.Range("L4").Value = 72.5
$ws.Range("M4").Value = -113.22223
$ws.Range("N4").Value = -300.5
$ws.Range("H12").Value = 28.125
$ws.Range("J12").Value = 15
$ws.Range("L12").Value = 15
$ws.Range("N12").Value = -355
$ws.Range("H116").Value = 7748.3125
$ws.Range("J116").Value = 8159.375
$ws.Range("L116").Value = 8159.375
$ws.Range("N116").Value = -15043.375
$ws.Range("H125").Value = 128067.25
$ws.Range("I125").Value = 252059.5
$ws.Range("J125").Value = 4075
$ws.Range("K125").Value = 2268535.5
$ws.Range("L125").Value = 36675
$ws.Range("M125").Value = -2266075.5
$ws.Range("N125").Value = -41595
$ws.Range("H132").Value = 1669.5588
$ws.Range("I132").Value = 1813.3572
$ws.Range("K132").Value = 5440.071599999999
$ws.Range("M132").Value = -2910.071599999999
$ws.Range("H133").Value = 77229.914
$ws.Range("J133").Value = 77229.914
$ws.Range("L133").Value = 77229.914
$ws.Range("N133").Value = -87349.914
$ws.Range("H134").Value = 99999
$ws.Range("J134").Value = 99999
$ws.Range("L134").Value = 99999
$ws.Range("N134").Value = -110139
$ws.Range("H136").Value = 96495.836
$ws.Range("J136").Value = 96495.836
$ws.Range("L136").Value = 96495.836
$ws.Range("N136").Value = -106695.836
$ws.Range("H137").Value = 317576
$ws.Range("J137").Value = 765931.9399999999
$ws.Range("L137").Value = 2297795.82
$ws.Range("N137").Value = -2302895.82
$ws.Range("H139").Value = 99999
$ws.Range("J139").Value = 99999
$ws.Range("L139").Value = 99999
$ws.Range("N139").Value = -110279
$ws.Range("H140").Value = 81794.375
$ws.Range("J140").Value = 81794.375
$ws.Range("L140").Value = 81794.375
$ws.Range("N140").Value = -92154.375

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9967.871999999999
$ws.Range("I32").Value = 5688.091
$ws.Range("K32").Value = 5688.091
$ws.Range("M32").Value = -5401.091
$ws.Range("H34").Value = 228341.67
$ws.Range("I34").Value = 35025
$ws.Range("J34").Value = 325000
$ws.Range("K34").Value = 35025
$ws.Range("L34").Value = 325000
$ws.Range("M34").Value = -34754
$ws.Range("N34").Value = -325542
$ws.Range("H45").Value = 3470.375
$ws.Range("I45").Value = 3772.4
$ws.Range("J45").Value = 2967
$ws.Range("K45").Value = 3772.4
$ws.Range("L45").Value = 2967
$ws.Range("M45").Value = -3395.4
$ws.Range("N45").Value = -3721
$ws.Range("H76").Value = 122571
$ws.Range("J76").Value = 122571
$ws.Range("L76").Value = 122571
$ws.Range("N76").Value = -123247
$ws.Range("H79").Value = 122571
$ws.Range("J79").Value = 122571
$ws.Range("L79").Value = 122571
$ws.Range("N79").Value = -124911
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H122").Value = 1877.2413
$ws.Range("I122").Value = 1662.1154
$ws.Range("J122").Value = 3741.6667
$ws.Range("K122").Value = 4986.3462
$ws.Range("L122").Value = 11225.0001
$ws.Range("M122").Value = -2536.3462
$ws.Range("N122").Value = -16125.0001
$ws.Range("H132").Value = 1750.25
$ws.Range("I132").Value = 1514.1892
$ws.Range("K132").Value = 4542.5676
$ws.Range("M132").Value = -2012.5676

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 73621.5
$ws.Range("J132").Value = 73621.5
$ws.Range("L132").Value = 73621.5
$ws.Range("N132").Value = -83741.5
$ws.Range("H134").Value = 3027.6667
$ws.Range("I134").Value = 918.5833
$ws.Range("K134").Value = 2755.7499
$ws.Range("M134").Value = -220.7498999999998
$ws.Range("H138").Value = 99999
$ws.Range("J138").Value = 99999
$ws.Range("L138").Value = 99999
$ws.Range("N138").Value = -110279
$ws.Range("H140").Value = 58166.484
$ws.Range("J140").Value = 43481.215
$ws.Range("L140").Value = 43481.215
$ws.Range("N140").Value = -53841.215

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 5714.4326
$ws.Range("I7").Value = 6449.5
$ws.Range("K7").Value = 6449.5
$ws.Range("M7").Value = -6336.5
$ws.Range("H22").Value = 675.7273
$ws.Range("I22").Value = 732.7778
$ws.Range("K22").Value = 732.7778
$ws.Range("M22").Value = -382.7778
$ws.Range("H31").Value = 2943.7058
$ws.Range("I31").Value = 2061.9167
$ws.Range("K31").Value = 2061.9167
$ws.Range("M31").Value = -1766.9167
$ws.Range("H34").Value = 2943.7058
$ws.Range("I34").Value = 2061.9167
$ws.Range("K34").Value = 2061.9167
$ws.Range("M34").Value = -1859.9167
$ws.Range("H132").Value = 2035.5714
$ws.Range("I132").Value = 2022.5454
$ws.Range("K132").Value = 6067.6362
$ws.Range("M132").Value = -3537.6362
$ws.Range("H138").Value = 108332.5
$ws.Range("J138").Value = 108332.5
$ws.Range("L138").Value = 108332.5
$ws.Range("N138").Value = -118612.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1035.4286
$ws.Range("J86").Value = 1141.3334
$ws.Range("L86").Value = 3424.0002
$ws.Range("N86").Value = -5796.0002
$ws.Range("H89").Value = 1035.4286
$ws.Range("J89").Value = 1141.3334
$ws.Range("L89").Value = 10272.0006
$ws.Range("N89").Value = -22128.0006
$ws.Range("H92").Value = 184.5
$ws.Range("I92").Value = 184.5
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 553.5
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 694.5
$ws.Range("N92").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5141.875
$ws.Range("I70").Value = 4937
$ws.Range("K70").Value = 4937
$ws.Range("M70").Value = -4667
$ws.Range("H73").Value = 5141.875
$ws.Range("I73").Value = 4937
$ws.Range("K73").Value = 4937
$ws.Range("M73").Value = -4001
$ws.Range("H92").Value = 9187.75
$ws.Range("J92").Value = 9187.75
$ws.Range("L92").Value = 9187.75
$ws.Range("N92").Value = -12931.75
$ws.Range("H107").Value = 1005.6111
$ws.Range("J107").Value = 1036.4546
$ws.Range("L107").Value = 1036.4546
$ws.Range("N107").Value = -4876.4546
$ws.Range("H126").Value = 4129
$ws.Range("I126").Value = 2712.5386
$ws.Range("J126").Value = 6430.75
$ws.Range("K126").Value = 8137.6158
$ws.Range("L126").Value = 19292.25
$ws.Range("M126").Value = -5667.6158
$ws.Range("N126").Value = -24232.25
$ws.Range("H135").Value = 52279.19
$ws.Range("J135").Value = 52279.19
$ws.Range("L135").Value = 52279.19
$ws.Range("N135").Value = -62419.19
$ws.Range("H140").Value = 97331.11
$ws.Range("J140").Value = 98185
$ws.Range("L140").Value = 98185
$ws.Range("N140").Value = -108545

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5188392
$ws.Range("I40").Value = 3561.9048
$ws.Range("J40").Value = 17286330
$ws.Range("K40").Value = 3561.9048
$ws.Range("L40").Value = 17286330
$ws.Range("M40").Value = -3425.9048
$ws.Range("N40").Value = -17286602
$ws.Range("H46").Value = 6131.037
$ws.Range("I46").Value = 9544.833000000001
$ws.Range("K46").Value = 9544.833000000001
$ws.Range("M46").Value = -9356.833000000001
$ws.Range("H55").Value = 7233.5
$ws.Range("I55").Value = 587.38464
$ws.Range("J55").Value = 36033.332
$ws.Range("K55").Value = 587.38464
$ws.Range("L55").Value = 36033.332
$ws.Range("M55").Value = -414.38464
$ws.Range("N55").Value = -36379.332
$ws.Range("H122").Value = 100004530
$ws.Range("I122").Value = 200004240
$ws.Range("K122").Value = 600012720
$ws.Range("M122").Value = -600010270
$ws.Range("H136").Value = 23258508
$ws.Range("I136").Value = 47622148
$ws.Range("K136").Value = 142866444
$ws.Range("M136").Value = -142863894

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 15000050
$ws.Range("I2").Value = 30000000
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 30000000
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = -29999888
$ws.Range("N2").Value = -324
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H26").Value = 1506
$ws.Range("I26").Value = 1012
$ws.Range("J26").Value = 2000
$ws.Range("K26").Value = 1012
$ws.Range("L26").Value = 2000
$ws.Range("M26").Value = -719
$ws.Range("N26").Value = -2586
$ws.Range("H107").Value = 1825.3572
$ws.Range("I107").Value = 356.57144
$ws.Range("J107").Value = 3294.1428
$ws.Range("K107").Value = 1069.71432
$ws.Range("L107").Value = 9882.428400000001
$ws.Range("M107").Value = 850.28568
$ws.Range("N107").Value = -13722.4284
$ws.Range("H122").Value = 2804.2222
$ws.Range("I122").Value = 2467.6
$ws.Range("K122").Value = 7402.799999999999
$ws.Range("M122").Value = -4952.799999999999
$ws.Range("H136").Value = 1460.1936
$ws.Range("I136").Value = 1257.125
$ws.Range("K136").Value = 3771.375
$ws.Range("M136").Value = -1221.375
